# The table currently holds one row per year from 2008 to 2020 (rows 2-14).
# The edit drops the 2008 and 2009 year rows and appends a new 2021 year row
# at the bottom, i.e. every remaining existing row shifts up by two and a
# fresh row of data is written for 2021.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete 2008-year and 2009-year rows; Excel shifts every row
# below them up accordingly.
$ws.Range("A2:A3").EntireRow.Delete()

# Write the new 2021-year row into what is now the last row (13). Row 13 is
# brand-new territory (beyond the previous last row), so first copy the
# year-label formatting (bold, centered, bordered) from the row above it,
# then fill in the values.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 0.6
$ws.Range("C13").Value = 12.7
$ws.Range("D13").Value = 70.9
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = 158.5
$ws.Range("G13").Value = 801.7
$ws.Range("H13").Value = 366.2
$ws.Range("I13").Value = 212.9
$ws.Range("J13").Value = 729.4
$ws.Range("K13").Value = 461.2
$ws.Range("L13").Value = 609.2
$ws.Range("M13").Value = ""
$ws.Range("N13").Value = 5.3
$ws.Range("O13").Value = 55.2
$ws.Range("P13").Value = 72.9
$ws.Range("Q13").Value = 1070.4
$ws.Range("R13").Value = 2.8
$ws.Range("S13").Value = 20.7
